# Update New Orleans xlsx: reorder sheets (review_info first, hotel_info
# second) and add a "State" column to hotel_info (inserted right after
# Hotel_Name, before City), populated with "Louisiana" for the existing row.

$wb = $excel.ActiveWorkbook

$wsHotel  = $wb.Worksheets.Item("hotel_info")
$wsReview = $wb.Worksheets.Item("review_info")

# --- Insert the new "State" column into hotel_info ------------------------
# Current layout: A=STR B=Hotel_Name C=City D=Zip E=TA_ReviewURL
#                 F=Tripadvisor_Hotel_Name G=English_Reviews_num
#                 H=Local_Rank I=Total_Reviews_num
# Insert a blank column at C so City (and everything after it) shifts right,
# then fill in the new column's header + value.
$wsHotel.Columns.Item(3).Insert()
$wsHotel.Range("C1").Value = "State"
$wsHotel.Range("C2").Value = "Louisiana"

# --- Reorder the sheets: review_info first, hotel_info second -------------
$wsReview.Move($wb.Worksheets.Item(1))
